# Applies the commit "Cai dat handlebar, scss, boostrapt":
#   - Remove the stray "DanhMuc"/numPr/ind paragraph formatting (w:pPr) from the
#     last three paragraphs of the "Thong tin du an" bullet list so they go back
#     to default (Normal) paragraph formatting.
#   - Bold + enlarge (sz 32) the "npm install express" run.
#   - Split "- cai dat nodemon ..." into three runs, capitalizing the leading "c".
#   - Fill in the previously-empty last paragraph with the morgan bullet text.

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph 9: "- Su dung Nodejs & Exprees : npm init de khoi tao, npm install express de cai express"
# (strip pPr, and bold+enlarge the "npm install express" run)
$p9 = $d.Paragraphs.Item(9)
$xml9 = $pkgOpen + '<w:p w14:paraId="3A10C8ED" w14:textId="23900AFF" w:rsidR="00AD4D6B" w:rsidRDefault="00AD4D6B" w:rsidP="00E0257C">' `
  + '<w:r><w:t xml:space="preserve">- Sử dụng Nodejs &amp; Exprees </w:t></w:r>' `
  + '<w:r w:rsidR="00E0257C"><w:t xml:space="preserve">: npm init để khởi tạo, </w:t></w:r>' `
  + '<w:r w:rsidR="00E0257C" w:rsidRPr="00E0257C"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>npm install express</w:t></w:r>' `
  + '<w:r w:rsidR="00E0257C"><w:t xml:space="preserve"> để cài express</w:t></w:r>' `
  + '</w:p>' + $pkgClose
$p9.Range.InsertXML($xml9)

# --- Paragraphs 10 & 11 are replaced together in a single InsertXML call: paragraph 11
# is the very last paragraph in the document body, and replacing only its own range
# leaves the old paragraph mark behind as a stray duplicate. Spanning the range across
# both paragraphs avoids that and replaces them cleanly.
#   10: "- cai dat nodemon de lang nghe su thay doi cua code" (capital C, split into runs)
#   11: previously empty, now holds the morgan bullet
$p10 = $d.Paragraphs.Item(10)
$p11 = $d.Paragraphs.Item(11)
$combined = $d.Range($p10.Range.Start, $p11.Range.End)
$xml1011 = $pkgOpen `
  + '<w:p w14:paraId="464734E1" w14:textId="3C014CAB" w:rsidR="00E0257C" w:rsidRDefault="00E0257C" w:rsidP="00E0257C">' `
  + '<w:r><w:t xml:space="preserve">- </w:t></w:r>' `
  + '<w:r><w:t>C</w:t></w:r>' `
  + '<w:r><w:t>ài đặt nodemon để lắng nghe sự thay đổi của code</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p w14:paraId="29AA7128" w14:textId="77777777" w:rsidR="00AD4D6B" w:rsidRPr="00AD4D6B" w:rsidRDefault="00AD4D6B" w:rsidP="00AD4D6B">' `
  + '<w:r><w:t>- Cài đặt morgan để lắng nghe các log từ client lên server</w:t></w:r>' `
  + '</w:p>' `
  + $pkgClose
$combined.InsertXML($xml1011)
